$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# ---------------------------------------------------------------------
# Week 1 (01-02 Dec) — fill in the journal entries that were left blank
# ---------------------------------------------------------------------

# Day 1 (A6, 01.12) block
$ws.Range("B6").Value = "Lecture du cahier des charges"
$ws.Range("D6").Value = 2

$ws.Range("B7").Value = "Création d'un business case"
$ws.Range("D7").Value = 1.5

# Day "rencontre" (A10, 01.12) block
$ws.Range("B11").Value = "Aide à Axelle pour le planning"
$ws.Range("D11").Value = 1.5

$ws.Range("B12").Value = "Users stories et product backlog"
$ws.Range("D12").Value = 2

# Day 2 (A14, 02.12) block
$ws.Range("B14").Value = "Préparation de la rencontre avec client"
$ws.Range("D14").Value = 2

# Row 15 becomes a highlighted "Rencontre" entry: unmerge B15:C15 and
# paint it with the same format as the orange "Rencontre" legend swatch.
$ws.Range("B15:C15").UnMerge()
$ws.Range("H5").Copy()
$ws.Range("B15:C15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Rencontre avec le client (kick-off)"
$ws.Range("D15").Value = 0.5

$ws.Range("B16").Value = "Documentation initiale"
$ws.Range("D16").Value = 1

# Weekly personal reflection for week 1
$ws.Range("A19").Value = "C'était la semaine de lancement du projet donc je n'ai pas forcément fait énormément de travail mais j'ai surtout réfléchi comment mener à bien le projet de la meilleure manière."

# ---------------------------------------------------------------------
# Week 3 (15-16 Dec) — fill in the journal entries that were left blank
# ---------------------------------------------------------------------

$ws.Range("B34").Value = "Documentation de la conception"
$ws.Range("D34").Value = 1

$ws.Range("B35").Value = "Configuration de la Raspberry"
$ws.Range("D35").Value = 1.5

$ws.Range("B36").Value = "Installation du serveur web sur la raspberry"
$ws.Range("D36").Value = 1

$ws.Range("B38").Value = "Déploiement sur la raspberry du PACMAN"
$ws.Range("D38").Value = 1.5

# Row 39 becomes a highlighted "Problème" entry: unmerge B39:C39 and
# paint it with the same format as the red "Problème" legend swatch.
$ws.Range("B39:C39").UnMerge()
$ws.Range("H2").Copy()
$ws.Range("B39:C39").PasteSpecial(-4122)
$ws.Range("B39").Value = "Test du jeu déployé (Bug de fluidité)"
$ws.Range("D39").Value = 0.5

# Row 40 becomes a highlighted "Solution" entry: unmerge B40:C40 and
# paint it with the same format as the green "Solution" legend swatch.
$ws.Range("B40:C40").UnMerge()
$ws.Range("H3").Copy()
$ws.Range("B40:C40").PasteSpecial(-4122)
$ws.Range("B40").Value = "Changement d'OS sur la raspberry pour optimiser la fluidité"
$ws.Range("D40").Value = 1.5

# ---------------------------------------------------------------------
# Update the selection state to match where the author ended up
# ---------------------------------------------------------------------
$null = $ws.Range("A19:D19").Select()
